$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.911.22"
$ws.Range("E2").Value = "  +0.26%  "

$ws.Range("D3").Value = "1.894.30"
$ws.Range("E3").Value = "  -0.03%  "

$ws.Range("E4").Value = "  +0.26%  "

$ws.Range("D5").Value = "'0.7741"
$ws.Range("E5").Value = "  -3.26%  "

$ws.Range("D6").Value = "'243.70"
$ws.Range("E6").Value = "  +0.40%  "

$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("E8").Value = "  -1.75%  "

$ws.Range("D9").Value = "'25.65"
$ws.Range("E9").Value = "  +0.13%  "

$ws.Range("D10").Value = "'0.07333"
$ws.Range("E10").Value = "  +3.98%  "

$ws.Range("D11").Value = "'0.08076"
$ws.Range("E11").Value = "  +0.33%  "

$ws.Range("D12").Value = "'0.7708"
$ws.Range("E12").Value = "  +0.15%  "

$ws.Range("D13").Value = "'5.494"
$ws.Range("E13").Value = "  +3.45%  "

$ws.Range("D14").Value = "1.880.27"
$ws.Range("E14").Value = "  -0.62%  "

$ws.Range("D15").Value = "'94.03"
$ws.Range("E15").Value = "  +1.83%  "

$ws.Range("D16").Value = "'6.221"
$ws.Range("E16").Value = "  +4.52%  "

$ws.Range("D17").Value = "29.912.43"
$ws.Range("E17").Value = "  +0.31%  "

$ws.Range("D18").Value = "'13.97"
$ws.Range("E18").Value = "  +0.68%  "

$ws.Range("D19").Value = "'247.20"
$ws.Range("E19").Value = "  +1.23%  "

$ws.Range("D20").Value = "'0.000007810"
$ws.Range("E20").Value = "  +1.23%  "

$ws.Range("D21").Value = "2.142.77"
$ws.Range("E21").Value = "  +0.51%  "

$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  +0.14%  "

$ws.Range("D23").Value = "'8.104"
$ws.Range("E23").Value = "  -1.79%  "

$ws.Range("E24").Value = "  +0.26%  "

$ws.Range("D25").Value = "'0.1580"
$ws.Range("E25").Value = "  -5.59%  "

$ws.Range("D26").Value = "'9.455"
$ws.Range("E26").Value = "  +1.20%  "

$ws.Range("D27").Value = "'163.18"
$ws.Range("E27").Value = "  -1.68%  "

$ws.Range("E28").Value = "  -0.07%  "

$ws.Range("E29").Value = "  -1.75%  "

$ws.Range("D30").Value = "'1.433"
$ws.Range("E30").Value = "  +2.75%  "

$ws.Range("D31").Value = "'1.543"
$ws.Range("E31").Value = "  +0.50%  "

$ws.Range("D32").Value = "'4.470"
$ws.Range("E32").Value = "  +1.12%  "

$ws.Range("D33").Value = "'0.05554"
$ws.Range("E33").Value = "  -2.09%  "

$ws.Range("D34").Value = "'4.061"
$ws.Range("E34").Value = "  -0.03%  "

$ws.Range("D35").Value = "'1.240"
$ws.Range("E35").Value = "  -1.84%  "

$ws.Range("D36").Value = "'0.7521"
$ws.Range("E36").Value = "  +1.56%  "

$ws.Range("E37").Value = "  -0.04%  "

$ws.Range("D38").Value = "'2.682"
$ws.Range("E38").Value = "  +2.25%  "

$ws.Range("D39").Value = "'0.01930"
$ws.Range("E39").Value = "  +0.92%  "

$ws.Range("D40").Value = "'2.793"
$ws.Range("E40").Value = "  +0.46%  "

$ws.Range("E41").Value = "  +0.96%  "

$ws.Range("D42").Value = "'73.99"
$ws.Range("E42").Value = "  +1.73%  "

$ws.Range("D43").Value = "1.099.42"
$ws.Range("E43").Value = "  +6.13%  "

$ws.Range("D44").Value = "'5.970"
$ws.Range("E44").Value = "  +2.60%  "

$ws.Range("D45").Value = "'0.8509"
$ws.Range("E45").Value = "  +0.59%  "

$ws.Range("D46").Value = "'0.9998"
$ws.Range("E46").Value = "  +0.14%  "

$ws.Range("D47").Value = "'1.886"
$ws.Range("E47").Value = "  +0.57%  "

$ws.Range("D48").Value = "'102.19"
$ws.Range("E48").Value = "  -0.40%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.776"
$ws.Range("E49").Value = "  -1.85%  "

$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").Value = "'7.519"
$ws.Range("E50").Value = "  +1.07%  "

$ws.Range("D51").Value = "'2.984"
$ws.Range("E51").Value = "  +1.79%  "
